$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header row (row 1) values:
# Before: A1=bedrooms_1, B1=living_rooms_1, C1=bedrooms_2, D1=living_rooms_2, E1=kitchens_1, F1=kitchens_2
# After:  A1=kitchens_1, B1=bedrooms_1, C1=living_rooms_1, D1=living_rooms_2, E1=kitchens_2, F1=bedrooms_2
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "living_rooms_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "kitchens_2"
$ws.Range("F1").Value = "bedrooms_2"
